$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: column E label changes from "PDF" to "Nome_arquivos_anexo"
$ws.Range("E1").Value = "Nome_arquivos_anexo"
$ws.Range("F1").Value = "STATUS"

# Add new data rows (emails + attachment file names)
$ws.Range("A2").Value = "oversouls11@gmail.com"
$ws.Range("E2").Value = "123456.pdf;21372.pdf"

$ws.Range("A3").Value = "financeiro4@webcertificados.com.br"
$ws.Range("E3").Value = "308828.pdf"

# Update the selected cell to mirror the author's last selection in the sheet
$ws.Range("G10").Select()
